$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2440182
$ws.Range("J17").Value = 2440182
$ws.Range("L17").Value = 7320546
$ws.Range("N17").Value = -7320882

$ws.Range("H32").Value = 4842.909
$ws.Range("I32").Value = 5493.25
$ws.Range("J32").Value = 4471.2856
$ws.Range("K32").Value = 5493.25
$ws.Range("L32").Value = 4471.2856
$ws.Range("M32").Value = -5167.25
$ws.Range("N32").Value = -5123.2856

$ws.Range("H33").Value = 1693.4286
$ws.Range("J33").Value = 3001
$ws.Range("L33").Value = 3001
$ws.Range("N33").Value = -3459

$ws.Range("H41").Value = 242
$ws.Range("I41").Value = 223
$ws.Range("K41").Value = 223
$ws.Range("M41").Value = 217

$ws.Range("H53").Value = 560.5238000000001
$ws.Range("I53").Value = 687.4545000000001
$ws.Range("J53").Value = 420.9
$ws.Range("K53").Value = 687.4545000000001
$ws.Range("L53").Value = 420.9
$ws.Range("M53").Value = -50.45450000000005
$ws.Range("N53").Value = -1694.9

$ws.Range("H62").Value = 2379.5
$ws.Range("I62").Value = 2379.5
$ws.Range("K62").Value = 2379.5
$ws.Range("M62").Value = -1755.5

$ws.Range("H65").Value = 2379.5
$ws.Range("I65").Value = 2379.5
$ws.Range("K65").Value = 11897.5
$ws.Range("M65").Value = -8777.5

$ws.Range("H86").Value = 2865.7856
$ws.Range("I86").Value = 2592.6
$ws.Range("J86").Value = 3548.75
$ws.Range("K86").Value = 2592.6
$ws.Range("L86").Value = 3548.75
$ws.Range("M86").Value = -1469.6
$ws.Range("N86").Value = -5794.75

$ws.Range("H89").Value = 2865.7856
$ws.Range("I89").Value = 2592.6
$ws.Range("J89").Value = 3548.75
$ws.Range("K89").Value = 12963
$ws.Range("L89").Value = 17743.75
$ws.Range("M89").Value = -7347
$ws.Range("N89").Value = -28975.75

$ws.Range("H132").Value = 6879.45
$ws.Range("I132").Value = 1449.2142
$ws.Range("K132").Value = 4347.642599999999
$ws.Range("M132").Value = -1817.642599999999

$ws.Range("H137").Value = 1466.8
$ws.Range("I137").Value = 1282.5714
$ws.Range("J137").Value = 1896.6666
$ws.Range("K137").Value = 3847.7142
$ws.Range("L137").Value = 5689.9998
$ws.Range("M137").Value = -1297.7142
$ws.Range("N137").Value = -10789.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2138.0833
$ws.Range("I2").Value = 758.9474
$ws.Range("K2").Value = 758.9474
$ws.Range("M2").Value = -645.9474

$ws.Range("H5").Value = 200.125
$ws.Range("I5").Value = 88.59999999999999
$ws.Range("K5").Value = 88.59999999999999
$ws.Range("M5").Value = 23.40000000000001

$ws.Range("H32").Value = 1193.5647
$ws.Range("I32").Value = 1221.0488
$ws.Range("K32").Value = 1221.0488
$ws.Range("M32").Value = -934.0488

$ws.Range("H61").Value = 10508.25
$ws.Range("I61").Value = 10547.5
$ws.Range("J61").Value = 10416.667
$ws.Range("K61").Value = 10547.5
$ws.Range("L61").Value = 10416.667
$ws.Range("M61").Value = -10335.5
$ws.Range("N61").Value = -10840.667

$ws.Range("H63").Value = 6646
$ws.Range("J63").Value = 9292
$ws.Range("L63").Value = 9292
$ws.Range("N63").Value = -10664

$ws.Range("H66").Value = 6646
$ws.Range("J66").Value = 9292
$ws.Range("L66").Value = 46460
$ws.Range("N66").Value = -53324

$ws.Range("H74").Value = 4797.423
$ws.Range("I74").Value = 4798.846
$ws.Range("J74").Value = 4796
$ws.Range("K74").Value = 4798.846
$ws.Range("L74").Value = 4796
$ws.Range("M74").Value = -3924.846
$ws.Range("N74").Value = -6544

$ws.Range("H77").Value = 4797.423
$ws.Range("I77").Value = 4798.846
$ws.Range("J77").Value = 4796
$ws.Range("K77").Value = 23994.23
$ws.Range("L77").Value = 23980
$ws.Range("M77").Value = -19626.23
$ws.Range("N77").Value = -32716

$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").Value = ""

$ws.Range("H116").Value = 2138.0833
$ws.Range("I116").Value = 758.9474
$ws.Range("K116").Value = 758.9474
$ws.Range("M116").Value = 1535.0526

$ws.Range("H132").Value = 3952.5098
$ws.Range("I132").Value = 3383.3333
$ws.Range("K132").Value = 10149.9999
$ws.Range("M132").Value = -7619.999899999999

$ws.Range("H136").Value = 10508.25
$ws.Range("I136").Value = 10547.5
$ws.Range("J136").Value = 10416.667
$ws.Range("K136").Value = 31642.5
$ws.Range("L136").Value = 31250.001
$ws.Range("M136").Value = -29092.5
$ws.Range("N136").Value = -36350.001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2138.0833
$ws.Range("I3").Value = 758.9474
$ws.Range("K3").Value = 758.9474
$ws.Range("M3").Value = -644.9474

$ws.Range("H4").Value = 200.125
$ws.Range("I4").Value = 88.59999999999999
$ws.Range("K4").Value = 88.59999999999999
$ws.Range("M4").Value = 26.40000000000001

$ws.Range("H26").Value = 20549.25
$ws.Range("I26").Value = 20549.25
$ws.Range("K26").Value = 20549.25
$ws.Range("M26").Value = -20257.25

$ws.Range("H134").Value = 3529.6875
$ws.Range("I134").Value = 3626.7173
$ws.Range("K134").Value = 10880.1519
$ws.Range("M134").Value = -8345.151899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5422.727
$ws.Range("I16").Value = 3777.5
$ws.Range("J16").Value = 7397
$ws.Range("K16").Value = 3777.5
$ws.Range("L16").Value = 7397
$ws.Range("M16").Value = -3490.5
$ws.Range("N16").Value = -7971

$ws.Range("H25").Value = 9999.5
$ws.Range("I25").Value = 9999.5
$ws.Range("K25").Value = 9999.5
$ws.Range("M25").Value = -9825.5

$ws.Range("H94").Value = 993.875
$ws.Range("I94").Value = 853.3333
$ws.Range("K94").Value = 853.3333
$ws.Range("M94").Value = -402.3333

$ws.Range("H113").Value = 5422.727
$ws.Range("I113").Value = 3777.5
$ws.Range("J113").Value = 7397
$ws.Range("K113").Value = 3777.5
$ws.Range("L113").Value = 7397
$ws.Range("M113").Value = -1607.5
$ws.Range("N113").Value = -11737

$ws.Range("H132").Value = 2634.5557
$ws.Range("I132").Value = 2244.4285
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 6733.2855
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -4203.2855
$ws.Range("N132").Value = -17060

$ws.Range("H134").Value = 5602.44
$ws.Range("I134").Value = 4801.65
$ws.Range("K134").Value = 14404.95
$ws.Range("M134").Value = -11869.95

$ws.Range("H138").Value = 126941.25
$ws.Range("I138").Value = 21995
$ws.Range("J138").Value = 161923.33
$ws.Range("K138").Value = 21995
$ws.Range("L138").Value = 161923.33
$ws.Range("M138").Value = -16855
$ws.Range("N138").Value = -172203.33

$ws.Range("H141").Value = 31765.334
$ws.Range("I141").Value = 27648
$ws.Range("K141").Value = 27648
$ws.Range("M141").Value = -22468

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 39106828
$ws.Range("I4").Value = 35416336
$ws.Range("J4").Value = 100000000
$ws.Range("K4").Value = 106249008
$ws.Range("L4").Value = 300000000
$ws.Range("M4").Value = -106248896
$ws.Range("N4").Value = -300000224

$ws.Range("H7").Value = 244.875
$ws.Range("I7").Value = 396
$ws.Range("J7").Value = 154.2
$ws.Range("K7").Value = 1188
$ws.Range("L7").Value = 462.6
$ws.Range("M7").Value = -1076
$ws.Range("N7").Value = -686.5999999999999

$ws.Range("H23").Value = 4886.5713
$ws.Range("J23").Value = 14389.143
$ws.Range("L23").Value = 43167.429
$ws.Range("N23").Value = -43637.429

$ws.Range("H131").Value = 31254052
$ws.Range("I131").Value = 71429810
$ws.Range("J131").Value = 6245.4443
$ws.Range("K131").Value = 214289430
$ws.Range("L131").Value = 18736.3329
$ws.Range("M131").Value = -214284390
$ws.Range("N131").Value = -28816.3329

$ws.Range("H141").Value = 1367.3
$ws.Range("I141").Value = 1367.3
$ws.Range("K141").Value = 4101.9
$ws.Range("M141").Value = 1078.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5394.143
$ws.Range("I80").Value = 4634.4
$ws.Range("K80").Value = 4634.4
$ws.Range("M80").Value = -3636.4

$ws.Range("H83").Value = 5394.143
$ws.Range("I83").Value = 4634.4
$ws.Range("K83").Value = 23172
$ws.Range("M83").Value = -18180

$ws.Range("H105").Value = 50000
$ws.Range("J105").Value = 50000
$ws.Range("L105").Value = 50000
$ws.Range("N105").Value = -56988

$ws.Range("H132").Value = 3119.6155
$ws.Range("I132").Value = 1243.8572
$ws.Range("K132").Value = 3731.5716
$ws.Range("M132").Value = -1201.5716

$ws.Range("H134").Value = 74994.336
$ws.Range("J134").Value = 74994.336
$ws.Range("L134").Value = 224983.008
$ws.Range("N134").Value = -230053.008

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 9461.105
$ws.Range("I132").Value = 11082.63
$ws.Range("K132").Value = 33247.89
$ws.Range("M132").Value = -30717.89

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4112.5713
$ws.Range("I126").Value = 3866.25
$ws.Range("K126").Value = 11598.75
$ws.Range("M126").Value = -9128.75

$ws.Range("H132").Value = 2903.638
$ws.Range("I132").Value = 2591
$ws.Range("K132").Value = 7773
$ws.Range("M132").Value = -5243

$ws.Range("H136").Value = 4649.8823
$ws.Range("I136").Value = 3323.2
$ws.Range("K136").Value = 9969.599999999999
$ws.Range("M136").Value = -7419.599999999999

$ws.Range("H137").Value = 79999.5
$ws.Range("J137").Value = 79999.5
$ws.Range("L137").Value = 79999.5
$ws.Range("N137").Value = -90199.5
